# "Generate Report for Archive"
#
# The localization-status report was regenerated: the "Ready for handoff"
# status became "In Translation" for the single tracked file, on every
# sheet that surfaces it (Overview's per-language columns, and each
# language sheet's own Status column). The narrower text also shrinks the
# autofit width of the columns that used to need room for the long
# "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview!E2 ("zh-cn" column) and Overview!F2 ("de-de" column)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Status column (C2) on each language sheet
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Column widths shrink to match the new, shorter status text (was
# 17.2159881591797 "chars", shrinks to 13.4101845877511 "chars").
# ColumnWidth round-trips through whole-pixel storage (snapped to 1/6-char
# steps), so feed it the nearest value that re-serializes to the target.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
